$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Diff": add the row describing the license change found
# for the Newtonsoft.Json nuget package.
# ---------------------------------------------------------------
$wsDiff = $wb.Worksheets.Item("Diff")

$wsDiff.Range("A2").Value = "nuget"
$wsDiff.Range("B2").Value = "Newtonsoft.Json"
$wsDiff.Range("C2").Value = "LICENSE_CHANGED"
$wsDiff.Range("D2").Value = "13.0.4"
$wsDiff.Range("E2").Value = "MIT"
$wsDiff.Range("F2").Value = "13.0.0"
$wsDiff.Range("G2").Value = "UNKNOWN"

$wsDiff.Hyperlinks.Add($wsDiff.Range("H2"), "https://www.nuget.org/packages/Newtonsoft.Json/13.0.0")

# Widen the columns so the newly added content fits (mirrors the
# auto-fit Excel performs after such an edit).
$wsDiff.Range("B1").ColumnWidth = 15.166666666666666
$wsDiff.Range("C1").ColumnWidth = 17.333333333333332
$wsDiff.Range("G1").ColumnWidth = 10
$wsDiff.Range("H1").ColumnWidth = 53.166666666666664

# ---------------------------------------------------------------
# Sheet "CurrentDependencies": update the existing Newtonsoft.Json
# entry with its new version / license / license url.
# ---------------------------------------------------------------
$wsCur = $wb.Worksheets.Item("CurrentDependencies")

$wsCur.Range("A2").Value = "nuget"
$wsCur.Range("B2").Value = "Newtonsoft.Json"
$wsCur.Range("C2").Value = "13.0.0"
$wsCur.Range("D2").Value = "UNKNOWN"

$wsCur.Hyperlinks.Add($wsCur.Range("E2"), "https://www.nuget.org/packages/Newtonsoft.Json/13.0.0")

$wsCur.Range("D1").ColumnWidth = 10
$wsCur.Range("E1").ColumnWidth = 53.166666666666664
